# Auto-generated PowerShell COM-interop script to update the cryptos list.
# Applies the per-cell value updates described by the commit diff
# ("Updated cryptos list on Thu Nov 21 23:50:49 UTC 2024 with GitHub Actions").
#
# Several of the new "Price" values are numeric-looking strings
# (e.g. "257.36", "35.80", "0.0000248", "1.00") that must keep their exact
# original text -- including trailing zeros and the thousands-dot grouping
# used elsewhere in this column -- rather than being auto-coerced to a
# number by Excel's normal cell-entry behaviour. For just those cells we
# force Text format ("@") before writing the value so the result stays a
# text string, matching the source data; every other cell is a plain
# Range.Value assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.547.09'
$ws.Range("E2").Value = '  +4.65%  '
$ws.Range("D3").Value = '3.368.39'
$ws.Range("E3").Value = '  +9.56%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.36'
$ws.Range("E5").Value = '  +8.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '623.72'
$ws.Range("E6").Value = '  +2.69%  '
$ws.Range("E7").Value = '  +10.35%  '
$ws.Range("E8").Value = '  +2.35%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '3.364.11'
$ws.Range("E10").Value = '  +9.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.816'
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("D13").Value = '98.277.91'
$ws.Range("E13").Value = '  +4.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.80'
$ws.Range("E14").Value = '  +6.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000248'
$ws.Range("E15").Value = '  +3.45%  '
$ws.Range("D16").Value = '3.996.79'
$ws.Range("E16").Value = '  +9.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.51'
$ws.Range("E17").Value = '  +3.98%  '
$ws.Range("D18").Value = '3.367.04'
$ws.Range("E18").Value = '  +9.76%  '
$ws.Range("E19").Value = '  +3.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.92'
$ws.Range("E20").Value = '  +4.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '486.72'
$ws.Range("E21").Value = '  +10.31%  '
$ws.Range("E22").Value = '  +3.24%  '
$ws.Range("E23").Value = '  +10.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.28'
$ws.Range("E24").Value = '  +5.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.79'
$ws.Range("E25").Value = '  +5.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.22'
$ws.Range("E26").Value = '  +4.45%  '
$ws.Range("E27").Value = '  +1.99%  '
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.188'
$ws.Range("E31").Value = '  +5.46%  '
$ws.Range("E32").Value = '  +2.62%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  +4.27%  '
$ws.Range("E35").Value = '  +7.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '522.56'
$ws.Range("E36").Value = '  +8.22%  '
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.38'
$ws.Range("E38").Value = '  -0.73%  '
$ws.Range("E39").Value = '  +4.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.82'
$ws.Range("E40").Value = '  +3.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.451'
$ws.Range("E41").Value = '  +4.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.28'
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.74'
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.27'
$ws.Range("E44").Value = '  +6.87%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.783'
$ws.Range("E45").Value = '  +16.31%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.35'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("E48").Value = '  +6.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.55'
$ws.Range("E49").Value = '  +6.84%  '
$ws.Range("E50").Value = '  +6.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.42'
$ws.Range("E51").Value = '  +4.22%  '
